$wb = $excel.ActiveWorkbook

# --- Sheet "2025" (sheet1.xml) ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 48306.03421671776
$ws.Range("B2").Value = 69142.60623028062
$ws.Range("E2").Value = 153393.0996716316
$ws.Range("I2").Value = 368437.3937326
$ws.Range("M2").Value = 117162.36729175
$ws.Range("N2").Value = 42606.7698102724
$ws.Range("O2").Value = 69179.85095077046

# --- Sheet "2030" (sheet2.xml) ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 11272.08507472299
$ws.Range("E2").Value = 259832.3234357048
$ws.Range("I2").Value = 303336.6638662838
$ws.Range("M2").Value = 108726.69049759
$ws.Range("N2").Value = 61774.6704292406
$ws.Range("O2").Value = 51779.80353542881

# --- Sheet "2035" (sheet3.xml) ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("E2").Value = 203665.4934311435
$ws.Range("G2").Value = 36325.45083359783
$ws.Range("I2").Value = 181475.9119087656
$ws.Range("L2").Value = 48934.58355751802
$ws.Range("M2").Value = 59245.21767383911
$ws.Range("N2").Value = 29435.62058728274
$ws.Range("O2").Value = 32090.85709793116
